$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with new values (replaces zero-filled rows with real readings)
$ws.Range("A2").Value = 45117.50694444445
$ws.Range("B2").Value = 17.295
$ws.Range("C2").Value = 11.407
$ws.Range("D2").Value = 3.988
$ws.Range("E2").Value = 36.822
$ws.Range("F2").Value = 28.934
$ws.Range("G2").Value = 13.611
$ws.Range("H2").Value = 42.104
$ws.Range("I2").Value = 20.942
$ws.Range("J2").Value = 8.632
$ws.Range("K2").Value = 12.848
$ws.Range("L2").Value = 14.482
$ws.Range("M2").Value = 15.027
$ws.Range("N2").Value = 4.344
$ws.Range("O2").Value = 13.535
$ws.Range("P2").Value = 18.799
$ws.Range("Q2").Value = 11.717
$ws.Range("R2").Value = 3.386
$ws.Range("S2").Value = 2.226
$ws.Range("T2").Value = 198.649
$ws.Range("U2").Value = 37.584
$ws.Range("V2").Value = 12.493
$ws.Range("W2").Value = 24.537
$ws.Range("X2").Value = 12.435
$ws.Range("Y2").Value = 3.148
$ws.Range("Z2").Value = 21.511
$ws.Range("AA2").Value = 11.035
$ws.Range("AB2").Value = 10.064
$ws.Range("AC2").Value = 11.833
$ws.Range("AD2").Value = 15.025
$ws.Range("AE2").Value = 3.317
$ws.Range("AF2").Value = 37.614
$ws.Range("AG2").Value = 6.758
$ws.Range("AH2").Value = 15.619
$ws.Range("A3").Value = 45117.51388888889
$ws.Range("B3").Value = 10.089
$ws.Range("C3").Value = 6.855
$ws.Range("D3").Value = 1.6
$ws.Range("E3").Value = 21.783
$ws.Range("F3").Value = 17.092
$ws.Range("G3").Value = 7.94
$ws.Range("H3").Value = 31.966
$ws.Range("I3").Value = 12.216
$ws.Range("J3").Value = 5.154
$ws.Range("K3").Value = 7.449
$ws.Range("L3").Value = 8.695
$ws.Range("M3").Value = 9.026999999999999
$ws.Range("N3").Value = 2.538
$ws.Range("O3").Value = 7.895
$ws.Range("P3").Value = 11.031
$ws.Range("Q3").Value = 6.99
$ws.Range("R3").Value = 1.477
$ws.Range("S3").Value = 0.898
$ws.Range("T3").Value = 112.859
$ws.Range("U3").Value = 22.174
$ws.Range("V3").Value = 7.288
$ws.Range("W3").Value = 14.481
$ws.Range("X3").Value = 7.445
$ws.Range("Y3").Value = 1.702
$ws.Range("Z3").Value = 15.414
$ws.Range("AA3").Value = 6.437
$ws.Range("AB3").Value = 5.928
$ws.Range("AC3").Value = 6.943
$ws.Range("AD3").Value = 9.016
$ws.Range("AE3").Value = 1.247
$ws.Range("AF3").Value = 29.412
$ws.Range("AG3").Value = 3.92
$ws.Range("AH3").Value = 9.112
$ws.Range("A4").Value = 45117.52083333334
$ws.Range("B4").Value = 0.961
$ws.Range("C4").Value = 0.253
$ws.Range("D4").Value = 0.804
$ws.Range("E4").Value = 2.02
$ws.Range("F4").Value = 0.997
$ws.Range("G4").Value = 0.763
$ws.Range("H4").Value = 8.577
$ws.Range("I4").Value = 1.163
$ws.Range("J4").Value = 0.452
$ws.Range("K4").Value = 0.238
$ws.Range("L4").Value = 0.802
$ws.Range("M4").Value = 0.737
$ws.Range("N4").Value = 0.277
$ws.Range("O4").Value = 0.752
$ws.Range("P4").Value = 1.079
$ws.Range("Q4").Value = 0.968
$ws.Range("R4").Value = 0.907
$ws.Range("S4").Value = 0.374
$ws.Range("T4").Value = 4.257
$ws.Range("U4").Value = 2.579
$ws.Range("V4").Value = 0.694
$ws.Range("W4").Value = 1.568
$ws.Range("X4").Value = 0.646
$ws.Range("Y4").Value = 0.541
$ws.Range("Z4").Value = 3.71
$ws.Range("AA4").Value = 0.613
$ws.Range("AB4").Value = 0.753
$ws.Range("AC4").Value = 0.841
$ws.Range("AD4").Value = 0.766
$ws.Range("AE4").Value = 0.766
$ws.Range("AF4").Value = 8.382999999999999
$ws.Range("AG4").Value = 0.222
$ws.Range("AH4").Value = 0.883
$ws.Range("A5").Value = 45117.52777777778
$ws.Range("B5").Value = 4.8
$ws.Range("C5").Value = 3.26
$ws.Range("D5").Value = 0.72
$ws.Range("E5").Value = 10.41
$ws.Range("F5").Value = 8.08
$ws.Range("G5").Value = 3.78
$ws.Range("H5").Value = 13.47
$ws.Range("I5").Value = 5.82
$ws.Range("J5").Value = 2.4
$ws.Range("K5").Value = 3.48
$ws.Range("L5").Value = 4.18
$ws.Range("M5").Value = 4.34
$ws.Range("N5").Value = 1.2
$ws.Range("O5").Value = 3.76
$ws.Range("P5").Value = 5.16
$ws.Range("Q5").Value = 3.39
$ws.Range("R5").Value = 0.71
$ws.Range("S5").Value = 0.38
$ws.Range("T5").Value = 49.84
$ws.Range("U5").Value = 10.35
$ws.Range("V5").Value = 3.47
$ws.Range("W5").Value = 6.66
$ws.Range("X5").Value = 3.52
$ws.Range("Y5").Value = 0.82
$ws.Range("Z5").Value = 6.31
$ws.Range("AA5").Value = 3.07
$ws.Range("AB5").Value = 2.84
$ws.Range("AC5").Value = 3.32
$ws.Range("AD5").Value = 4.34
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 12.03
$ws.Range("AG5").Value = 1.84
$ws.Range("AH5").Value = 4.34

# Row 6 no longer present in the new dataset
$ws.Rows.Item(6).Delete()

# New column widths for data columns B:AH
$ws.Range("B1").ColumnWidth = 7.166666666666667
$ws.Range("C1").ColumnWidth = 7.166666666666667
$ws.Range("D1").ColumnWidth = 6.166666666666667
$ws.Range("E1").ColumnWidth = 7.166666666666667
$ws.Range("F1").ColumnWidth = 7.166666666666667
$ws.Range("G1").ColumnWidth = 7.166666666666667
$ws.Range("H1").ColumnWidth = 7.166666666666667
$ws.Range("I1").ColumnWidth = 7.166666666666667
$ws.Range("J1").ColumnWidth = 6.166666666666667
$ws.Range("K1").ColumnWidth = 7.166666666666667
$ws.Range("L1").ColumnWidth = 7.166666666666667
$ws.Range("M1").ColumnWidth = 7.166666666666667
$ws.Range("N1").ColumnWidth = 6.166666666666667
$ws.Range("O1").ColumnWidth = 7.166666666666667
$ws.Range("P1").ColumnWidth = 7.166666666666667
$ws.Range("Q1").ColumnWidth = 7.166666666666667
$ws.Range("R1").ColumnWidth = 6.166666666666667
$ws.Range("S1").ColumnWidth = 6.166666666666667
$ws.Range("T1").ColumnWidth = 8.166666666666666
$ws.Range("U1").ColumnWidth = 7.166666666666667
$ws.Range("V1").ColumnWidth = 7.166666666666667
$ws.Range("W1").ColumnWidth = 7.166666666666667
$ws.Range("X1").ColumnWidth = 7.166666666666667
$ws.Range("Y1").ColumnWidth = 6.166666666666667
$ws.Range("Z1").ColumnWidth = 7.166666666666667
$ws.Range("AA1").ColumnWidth = 7.166666666666667
$ws.Range("AB1").ColumnWidth = 7.166666666666667
$ws.Range("AC1").ColumnWidth = 7.166666666666667
$ws.Range("AD1").ColumnWidth = 7.166666666666667
$ws.Range("AE1").ColumnWidth = 6.166666666666667
$ws.Range("AF1").ColumnWidth = 7.166666666666667
$ws.Range("AG1").ColumnWidth = 6.166666666666667
$ws.Range("AH1").ColumnWidth = 7.166666666666667
